$wb = $excel.ActiveWorkbook

# Both worksheets ("dust_event" and "dust_event_small") share the same
# File-Location/Storage lookup table in columns A:B. Update the lookup
# labels and drop the now-unused rows (5-8) from that table on both sheets.
foreach ($ws in $wb.Worksheets) {
    $ws.Range("B3").Value = "Processing"
    $ws.Range("B4").Value = "Unknown"
    $ws.Range("A5:B8").ClearContents()
}

# Flip which sheet is active/selected: "dust_event_small" (2nd sheet)
# becomes the active tab, with a new selection on each sheet.
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Range("C5").Select() | Out-Null
$ws2.Select() | Out-Null
$ws2.Range("B10").Select() | Out-Null
